$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Fix product-name text: add missing hyphen after "245"
$productName = "245-MS-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"
$wsInput.Range("B1").Value = $productName
$wsOutput.Range("B1").Value = $productName

# Move the selection to B1 on the input sheet, then switch focus/selection
# to the output sheet (also at B1), matching the saved view state.
$wsInput.Range("B1").Select()
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
